# Refresh the cryptos price list (Price / Volume(1h) columns, plus a
# handful of rank swaps that moved coin name/link/price together) to
# match the latest GitHub Actions scrape.
#
# Cells in column D that now hold a plain decimal (e.g. "1.00", "0.996")
# are written with a leading apostrophe so Excel keeps them as text
# (matching the sheet's existing inlineStr/text cells) instead of
# silently re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.730.47"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "2.692.86"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'525.28"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "'145.32"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +1.51%  "
$ws.Range("D9").Value = "2.711.28"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'6.51"
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "3.163.14"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").Value = "60.684.40"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "'21.36"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.737.73"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000139"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "'349.70"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "'4.52"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "'10.58"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("D22").Value = "'6.33"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'63.65"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").Value = "'0.422"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.0₃0822"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").Value = "'7.34"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").Value = "'6.90"
$ws.Range("E30").Value = "  +8.99%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'19.26"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").Value = "'150.07"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "'4.27"
$ws.Range("E35").Value = "  +6.08%  "
$ws.Range("E36").Value = "  +10.18%  "
$ws.Range("D37").Value = "'0.950"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("D38").Value = "'0.881"
$ws.Range("E38").Value = "  +4.34%  "
$ws.Range("E39").Value = "  +8.15%  "
$ws.Range("D40").Value = "'36.96"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").Value = "'283.40"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "'20.14"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.612"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.144.14"
$ws.Range("E46").Value = "  +7.31%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'0.995"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  +4.86%  "
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").Value = "'10.47"
$ws.Range("E51").Value = "  +1.86%  "
